$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Q3 and R3 with rounded values
$ws.Range("Q3").Value = 575010
$ws.Range("R3").Value = 6299808

# Clear the Starttid (Z3) and Sluttid (AB3) cells entirely
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
